$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.864.85'
$ws.Range('E2').Value = '  +0.57%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.767.67'
$ws.Range('E3').Value = '  +0.54%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.59'
$ws.Range('E5').Value = '  +0.88%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.07%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4479'
$ws.Range('E7').Value = '  -2.58%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3558'
$ws.Range('E8').Value = '  -1.26%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07428'
$ws.Range('E9').Value = '  -1.29%  '

# Row 10
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.06'
$ws.Range('E10').Value = '  -0.26%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.095'
$ws.Range('E11').Value = '  -0.70%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.00%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.86'
$ws.Range('E13').Value = '  +0.06%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.035'
$ws.Range('E14').Value = '  +0.05%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.204'
$ws.Range('E15').Value = '  +1.07%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.768.16'
$ws.Range('E16').Value = '  +0.49%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.89'
$ws.Range('E17').Value = '  +0.58%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001057'
$ws.Range('E18').Value = '  -0.86%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06434'
$ws.Range('E19').Value = '  +0.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.08%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.18'
$ws.Range('E21').Value = '  +2.16%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.806'
$ws.Range('E22').Value = '  -0.46%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.897.71'
$ws.Range('E23').Value = '  +0.53%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.29'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.109'
$ws.Range('E25').Value = '  -0.16%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.61'
$ws.Range('E26').Value = '  -1.36%  '

# Row 27
$ws.Range('E27').Value = '  -1.08%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.977.69'
$ws.Range('E28').Value = '  +0.98%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.155'
$ws.Range('E29').Value = '  +2.55%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.35'
$ws.Range('E30').Value = '  -1.08%  '

# Row 31
$ws.Range('E31').Value = '  +2.88%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09139'
$ws.Range('E32').Value = '  -1.21%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.583'
$ws.Range('E33').Value = '  +0.41%  '

# Row 34
$ws.Range('E34').Value = '  -0.93%  '

# Row 35
$ws.Range('E35').Value = '  -1.14%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02294'
$ws.Range('E36').Value = '  -0.61%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06098'
$ws.Range('E37').Value = '  +0.61%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2095'
$ws.Range('E38').Value = '  -0.74%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6323'
$ws.Range('E39').Value = '  -0.96%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.963'
$ws.Range('E40').Value = '  -0.61%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.183'
$ws.Range('E41').Value = '  -1.65%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.390'
$ws.Range('E42').Value = '  +0.78%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.919'
$ws.Range('E43').Value = '  +0.61%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.24'
$ws.Range('E44').Value = '  -0.89%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.735'
$ws.Range('E45').Value = '  +0.54%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5876'
$ws.Range('E46').Value = '  -0.80%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.38'
$ws.Range('E47').Value = '  -0.90%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955'
$ws.Range('E48').Value = '  -0.34%  '

# Row 49
$ws.Range('E49').Value = '  +0.56%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.137'
$ws.Range('E50').Value = '  -1.78%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.86'
$ws.Range('E51').Value = '  +0.37%  '
